# Normalize the "保險" (insurance, sheet7) and "債務" (debt, sheet8) sheets so
# their header row + metadata columns match every other sheet in the workbook
# (property_category/category/date/legislator_name/legislator_id/source_file/index).

$wb = $excel.ActiveWorkbook

$legislatorName = "王進士"
$legislatorId = "1701"
$sourceFile = "tmpf41"
$category = "normal"
$theDate = "2011-12-28"

# ---------------------------------------------------------------------------
# Sheet7 "保險" (insurance)
# ---------------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(7)

# Copy the existing bold/bordered header style onto the new header cells
# (F1:K1) before writing their labels, so they match B1:D1's formatting.
$wsIns.Range("B1").Copy()
$wsIns.Range("E1:K1").PasteSpecial(-4122)

# Fix up the data rows first (values were duplicated/shifted before the fix).
$wsIns.Cells.Item(2,2).Value = "國泰人壽"
$wsIns.Cells.Item(2,3).Value = "超越變動額萬能壽險乙型"
$wsIns.Cells.Item(2,4).Value = "周麗容"
$wsIns.Cells.Item(2,5).Value = "insurance"
$wsIns.Cells.Item(2,6).Value = $category
$wsIns.Cells.Item(2,7).Value = $theDate
$wsIns.Cells.Item(2,8).Value = $legislatorName
$wsIns.Cells.Item(2,9).Value = 1701
$wsIns.Cells.Item(2,10).Value = $sourceFile
$wsIns.Cells.Item(2,11).Value = 117

$wsIns.Cells.Item(3,2).Value = "台灣銀行"
$wsIns.Cells.Item(3,3).Value = "美麗人生萬能保險"
$wsIns.Cells.Item(3,4).Value = "王進士"
$wsIns.Cells.Item(3,5).Value = "insurance"
$wsIns.Cells.Item(3,6).Value = $category
$wsIns.Cells.Item(3,7).Value = $theDate
$wsIns.Cells.Item(3,8).Value = $legislatorName
$wsIns.Cells.Item(3,9).Value = 1701
$wsIns.Cells.Item(3,10).Value = $sourceFile
$wsIns.Cells.Item(3,11).Value = 118

# Now write the real header labels on row 1 (replacing the stray leftover
# data that used to sit there), clearing the old E1 numeric value.
$wsIns.Cells.Item(1,2).Value = "company"
$wsIns.Cells.Item(1,3).Value = "name"
$wsIns.Cells.Item(1,4).Value = "owner"
$wsIns.Cells.Item(1,5).Value = "property_category"
$wsIns.Cells.Item(1,6).Value = "category"
$wsIns.Cells.Item(1,7).Value = "date"
$wsIns.Cells.Item(1,8).Value = "legislator_name"
$wsIns.Cells.Item(1,9).Value = "legislator_id"
$wsIns.Cells.Item(1,10).Value = "source_file"
$wsIns.Cells.Item(1,11).Value = "index"

# ---------------------------------------------------------------------------
# Sheet8 "債務" (debt)
# ---------------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(8)

# Copy the existing bold/bordered header style onto the new header cells
# (H1:N1) before writing their labels.
$wsDebt.Range("B1").Copy()
$wsDebt.Range("H1:N1").PasteSpecial(-4122)

# Data row 2: append the metadata columns (existing B2:G2 stay as-is).
$wsDebt.Cells.Item(2,8).Value = "debt"
$wsDebt.Cells.Item(2,9).Value = $category
$wsDebt.Cells.Item(2,10).Value = $theDate
$wsDebt.Cells.Item(2,11).Value = $legislatorName
$wsDebt.Cells.Item(2,12).Value = 1701
$wsDebt.Cells.Item(2,13).Value = $sourceFile
$wsDebt.Cells.Item(2,14).Value = 128

# Now write the real header labels on row 1.
$wsDebt.Cells.Item(1,2).Value = "species"
$wsDebt.Cells.Item(1,3).Value = "debtor"
$wsDebt.Cells.Item(1,4).Value = "owner"
$wsDebt.Cells.Item(1,5).Value = "total"
$wsDebt.Cells.Item(1,6).Value = "register_date"
$wsDebt.Cells.Item(1,7).Value = "register_reason"
$wsDebt.Cells.Item(1,8).Value = "property_category"
$wsDebt.Cells.Item(1,9).Value = "category"
$wsDebt.Cells.Item(1,10).Value = "date"
$wsDebt.Cells.Item(1,11).Value = "legislator_name"
$wsDebt.Cells.Item(1,12).Value = "legislator_id"
$wsDebt.Cells.Item(1,13).Value = "source_file"
$wsDebt.Cells.Item(1,14).Value = "index"
